# 导入导出忽略 tenant_id, tenant_id_lbl 等
# Remove the "tenant_id" column (header in K1: comment.tenant_id_lbl,
# data in K2: model.tenant_id_lbl) from the Sheet1 template, shifting the
# trailing "update_time" column (L -> K) left to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("K").Delete()
